# PowerEdges.xlsx edit: "compare the metrics and plot distribution"
#
# The underlying edge list in columns B (START POWER NODE ID) and C (END
# POWER NODE ID) gets normalized so that, for every row, the smaller of the
# two node ids sits in column B and the larger sits in column C (several
# rows previously had them the other way around). Rows that already
# satisfied B < C are left untouched.
#
# The workbook was also left scrolled near the bottom of the data with
# B77 (the first blank row right after the data, row 77 = header + 76 data
# rows) selected, so we reproduce that final selection too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 76

for ($r = 2; $r -le $lastRow; $r++) {
    $startCell = $ws.Cells.Item($r, 2)
    $endCell   = $ws.Cells.Item($r, 3)

    $startVal = $startCell.Value2
    $endVal   = $endCell.Value2

    if ($startVal -gt $endVal) {
        $startCell.Value2 = $endVal
        $endCell.Value2   = $startVal
    }
}

# Leave the sheet scrolled to the end of the list with B77 selected, matching
# where the author ended up after reviewing/fixing the data.
$ws.Range("B77").Select()
